# Add two new worksheets ("valid_login" and "invalid_login") after the
# existing "Sheet1", populate them with login test data, and give the new
# cells the same cell style the workbook's existing cells already use
# (so no new style/font entries are introduced).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Insert the new sheets right after Sheet1, in order.
$wsValid = $wb.Worksheets.Add([System.Type]::Missing, $ws1)
$wsValid.Name = "valid_login"

$wsInvalid = $wb.Worksheets.Add([System.Type]::Missing, $wsValid)
$wsInvalid.Name = "invalid_login"

# --- valid_login ---------------------------------------------------------
$wsValid.Range("A1").Value = "User_name"
$wsValid.Range("B1").Value = "Password"
$wsValid.Range("A2").Value = "admin"
$wsValid.Range("B2").Value = "manager"
$wsValid.Range("A3").Value = "trainee"
$wsValid.Range("B3").Value = "trainee"

# --- invalid_login ---------------------------------------------------------
$wsInvalid.Range("A1").Value = "User_name"
$wsInvalid.Range("B1").Value = "Password"
$wsInvalid.Range("A2").Value = "abcd"
$wsInvalid.Range("B2").Value = "xyz"

# Match the existing workbook's default cell style (same one already used
# on Sheet1) instead of leaving the brand-new cells on the generic style,
# by copying formatting only from an already-styled cell.
$ws1.Range("A1").Copy()
$wsValid.Range("A1:B3").PasteSpecial(-4122)
$wsInvalid.Range("A1:B2").PasteSpecial(-4122)

$excel.CutCopyMode = 0
